$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing the original text formatting,
# e.g. trailing zeros like "4.30" -> 4.3).
$textCells = @("D5", "D10", "D11", "D15", "D16", "D19", "D20", "D22", "D24", "D25", "D29", "D31", "D37", "D39", "D40", "D42", "D44", "D46", "D47", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.254.50"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.591.86"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "212.85"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").Value = "18.96"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("D11").Value = "0.0851"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "1.816.57"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.599.84"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "0.509"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").Value = "63.78"
$ws.Range("D17").Value = "26.260.63"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "215.59"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "7.34"
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "4.30"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").Value = "2.10"
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("D25").Value = "145.01"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").Value = "15.10"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").Value = "1.421.18"
$ws.Range("E33").Value = "  +6.02%  "
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").Value = "0.572"
$ws.Range("E37").Value = "  -4.32%  "
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").Value = "0.826"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").Value = "5.77"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "0.941"
$ws.Range("E42").Value = "  -10.62%  "
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").Value = "0.762"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "1.728.62"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").Value = "60.89"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("D47").Value = "86.84"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").Value = "0.0952"
$ws.Range("E50").Value = "  -2.84%  "
$ws.Range("E51").Value = "  +0.03%  "
